$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81; this shifts the existing rows 81-133 down to 82-134
$ws.Rows(81).Insert()

# Populate the newly inserted row 81 with the new DNS record
$ws.Range("A81").Value = "868b4cabc8c6415f29c34b34bfa72210"
$ws.Range("B81").Value = "bin.italiacdn.net"
$ws.Range("C81").Value = "CNAME"
$ws.Range("D81").Value = "62012bb4-9a36-4017-b533-489f1e455049.cfargotunnel.com"
$ws.Range("E81").Value = $true
$ws.Range("F81").Value = $true
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = "{'flatten_cname': False}"
$ws.Range("I81").Value = "{}"
$ws.Range("J81").Value = ""
$ws.Range("K81").Value = "[]"
$ws.Range("L81").Value = "2025-03-17T21:43:23.250682Z"
$ws.Range("M81").Value = "2025-03-17T21:43:23.250682Z"
$ws.Range("N81").Value = ""
$ws.Range("O81").Value = ""
